# Update countries & provincias Spain
# Refresh the COVID-19 "Pais" sheet with the latest data pull:
#  - bump the "Datos actualizados" timestamp
#  - update case/recovered/death counters for the countries that changed
#  - a handful of neighbouring countries swapped ranking order (same
#    totals moved to the other row), so their names are rewritten too
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 1
$ws.Cells.Item(1, 1).Value = "Datos actualizados a 9 de Julio de 2020 a las 10:34"

# Row 4
$ws.Cells.Item(4, 1).Value = "Estados Unidos"
$ws.Cells.Item(4, 2).Value = 3159414
$ws.Cells.Item(4, 3).Value = 482
$ws.Cells.Item(4, 4).Value = 1392808
$ws.Cells.Item(4, 5).Value = 1631739
$ws.Cells.Item(4, 6).Value = 0
$ws.Cells.Item(4, 7).Value = 5
$ws.Cells.Item(4, 8).Value = 134867

# Row 6
$ws.Cells.Item(6, 1).Value = "India"
$ws.Cells.Item(6, 2).Value = 769257
$ws.Cells.Item(6, 3).Value = 205
$ws.Cells.Item(6, 4).Value = 476600
$ws.Cells.Item(6, 5).Value = 271496
$ws.Cells.Item(6, 6).Value = 0
$ws.Cells.Item(6, 7).Value = 17
$ws.Cells.Item(6, 8).Value = 21161

# Row 7
$ws.Cells.Item(7, 1).Value = "Rusia"
$ws.Cells.Item(7, 2).Value = 707301
$ws.Cells.Item(7, 3).Value = 6509
$ws.Cells.Item(7, 4).Value = 481316
$ws.Cells.Item(7, 5).Value = 215142
$ws.Cells.Item(7, 6).Value = 0
$ws.Cells.Item(7, 7).Value = 176
$ws.Cells.Item(7, 8).Value = 10843

# Row 33
$ws.Cells.Item(33, 1).Value = "Belgica"
$ws.Cells.Item(33, 2).Value = 62210
$ws.Cells.Item(33, 3).Value = 87
$ws.Cells.Item(33, 4).Value = 17159
$ws.Cells.Item(33, 5).Value = 35273
$ws.Cells.Item(33, 6).Value = 0
$ws.Cells.Item(33, 7).Value = 2
$ws.Cells.Item(33, 8).Value = 9778

# Row 41
$ws.Cells.Item(41, 1).Value = "Singapur"
$ws.Cells.Item(41, 2).Value = 45423
$ws.Cells.Item(41, 3).Value = 125
$ws.Cells.Item(41, 4).Value = 41323
$ws.Cells.Item(41, 5).Value = 4074
$ws.Cells.Item(41, 6).Value = 0
$ws.Cells.Item(41, 7).Value = 0
$ws.Cells.Item(41, 8).Value = 26

# Row 46
$ws.Cells.Item(46, 1).Value = "Polonia"
$ws.Cells.Item(46, 2).Value = 36689
$ws.Cells.Item(46, 3).Value = 0
$ws.Cells.Item(46, 4).Value = 25477
$ws.Cells.Item(46, 5).Value = 9670
$ws.Cells.Item(46, 6).Value = 0
$ws.Cells.Item(46, 7).Value = 0
$ws.Cells.Item(46, 8).Value = 1542

# Row 47
$ws.Cells.Item(47, 1).Value = "Israel"
$ws.Cells.Item(47, 2).Value = 33947
$ws.Cells.Item(47, 3).Value = 390
$ws.Cells.Item(47, 4).Value = 18392
$ws.Cells.Item(47, 5).Value = 15209
$ws.Cells.Item(47, 6).Value = 0
$ws.Cells.Item(47, 7).Value = 2
$ws.Cells.Item(47, 8).Value = 346

# Row 48
$ws.Cells.Item(48, 1).Value = "Afganistan"
$ws.Cells.Item(48, 2).Value = 33908
$ws.Cells.Item(48, 3).Value = 314
$ws.Cells.Item(48, 4).Value = 20847
$ws.Cells.Item(48, 5).Value = 12104
$ws.Cells.Item(48, 6).Value = 0
$ws.Cells.Item(48, 7).Value = 21
$ws.Cells.Item(48, 8).Value = 957

# Row 50
$ws.Cells.Item(50, 1).Value = "Barein"
$ws.Cells.Item(50, 2).Value = 30931
$ws.Cells.Item(50, 3).Value = 0
$ws.Cells.Item(50, 4).Value = 26073
$ws.Cells.Item(50, 5).Value = 4757
$ws.Cells.Item(50, 6).Value = 0
$ws.Cells.Item(50, 7).Value = 3
$ws.Cells.Item(50, 8).Value = 101

# Row 76
$ws.Cells.Item(76, 1).Value = "El Salvador"
$ws.Cells.Item(76, 2).Value = 8844
$ws.Cells.Item(76, 3).Value = 278
$ws.Cells.Item(76, 4).Value = 5289
$ws.Cells.Item(76, 5).Value = 3312
$ws.Cells.Item(76, 6).Value = 0
$ws.Cells.Item(76, 7).Value = 8
$ws.Cells.Item(76, 8).Value = 243

# Row 77
$ws.Cells.Item(77, 1).Value = "Malasia"
$ws.Cells.Item(77, 2).Value = 8677
$ws.Cells.Item(77, 3).Value = 0
$ws.Cells.Item(77, 4).Value = 8486
$ws.Cells.Item(77, 5).Value = 70
$ws.Cells.Item(77, 6).Value = 0
$ws.Cells.Item(77, 7).Value = 0
$ws.Cells.Item(77, 8).Value = 121

# Row 97
$ws.Cells.Item(97, 1).Value = "Hungria"
$ws.Cells.Item(97, 2).Value = 4220
$ws.Cells.Item(97, 3).Value = 10
$ws.Cells.Item(97, 4).Value = 2887
$ws.Cells.Item(97, 5).Value = 742
$ws.Cells.Item(97, 6).Value = 0
$ws.Cells.Item(97, 7).Value = 2
$ws.Cells.Item(97, 8).Value = 591

# Row 114
$ws.Cells.Item(114, 1).Value = "Estonia"
$ws.Cells.Item(114, 2).Value = 2011
$ws.Cells.Item(114, 3).Value = 8
$ws.Cells.Item(114, 4).Value = 1889
$ws.Cells.Item(114, 5).Value = 53
$ws.Cells.Item(114, 6).Value = 0
$ws.Cells.Item(114, 7).Value = 0
$ws.Cells.Item(114, 8).Value = 69

# Row 120
$ws.Cells.Item(120, 1).Value = "Eslovaquia"
$ws.Cells.Item(120, 2).Value = 1851
$ws.Cells.Item(120, 3).Value = 53
$ws.Cells.Item(120, 4).Value = 1477
$ws.Cells.Item(120, 5).Value = 346
$ws.Cells.Item(120, 6).Value = 0
$ws.Cells.Item(120, 7).Value = 0
$ws.Cells.Item(120, 8).Value = 28

# Row 121
$ws.Cells.Item(121, 1).Value = "Congo"
$ws.Cells.Item(121, 2).Value = 1821
$ws.Cells.Item(121, 3).Value = 0
$ws.Cells.Item(121, 4).Value = 525
$ws.Cells.Item(121, 5).Value = 1249
$ws.Cells.Item(121, 6).Value = 0
$ws.Cells.Item(121, 7).Value = 0
$ws.Cells.Item(121, 8).Value = 47

# Row 123
$ws.Cells.Item(123, 1).Value = "Eslovenia"
$ws.Cells.Item(123, 2).Value = 1776
$ws.Cells.Item(123, 3).Value = 13
$ws.Cells.Item(123, 4).Value = 1429
$ws.Cells.Item(123, 5).Value = 236
$ws.Cells.Item(123, 6).Value = 0
$ws.Cells.Item(123, 7).Value = 0
$ws.Cells.Item(123, 8).Value = 111

# Row 134
$ws.Cells.Item(134, 1).Value = "Letonia"
$ws.Cells.Item(134, 2).Value = 1154
$ws.Cells.Item(134, 3).Value = 13
$ws.Cells.Item(134, 4).Value = 1019
$ws.Cells.Item(134, 5).Value = 105
$ws.Cells.Item(134, 6).Value = 0
$ws.Cells.Item(134, 7).Value = 0
$ws.Cells.Item(134, 8).Value = 30

# Row 184
$ws.Cells.Item(184, 1).Value = "Seychelles"
$ws.Cells.Item(184, 2).Value = 91
$ws.Cells.Item(184, 3).Value = 0
$ws.Cells.Item(184, 4).Value = 11
$ws.Cells.Item(184, 5).Value = 80
$ws.Cells.Item(184, 6).Value = 0
$ws.Cells.Item(184, 7).Value = 0
$ws.Cells.Item(184, 8).Value = 0

# Row 185
$ws.Cells.Item(185, 1).Value = "Lesoto"
$ws.Cells.Item(185, 2).Value = 91
$ws.Cells.Item(185, 3).Value = 0
$ws.Cells.Item(185, 4).Value = 11
$ws.Cells.Item(185, 5).Value = 80
$ws.Cells.Item(185, 6).Value = 0
$ws.Cells.Item(185, 7).Value = 0
$ws.Cells.Item(185, 8).Value = 0

# Row 209
$ws.Cells.Item(209, 1).Value = "Islas Malvinas"
$ws.Cells.Item(209, 2).Value = 13
$ws.Cells.Item(209, 3).Value = 0
$ws.Cells.Item(209, 4).Value = 13
$ws.Cells.Item(209, 5).Value = 0
$ws.Cells.Item(209, 6).Value = 0
$ws.Cells.Item(209, 7).Value = 0
$ws.Cells.Item(209, 8).Value = 0

# Row 210
$ws.Cells.Item(210, 1).Value = "Groenlandia"
$ws.Cells.Item(210, 2).Value = 13
$ws.Cells.Item(210, 3).Value = 0
$ws.Cells.Item(210, 4).Value = 13
$ws.Cells.Item(210, 5).Value = 0
$ws.Cells.Item(210, 6).Value = 0
$ws.Cells.Item(210, 7).Value = 0
$ws.Cells.Item(210, 8).Value = 0
